$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 629.6667
$ws.Range("I2").Value = 329
$ws.Range("J2").Value = 689.8
$ws.Range("K2").Value = 329
$ws.Range("L2").Value = 689.8
$ws.Range("M2").Value = -216
$ws.Range("N2").Value = -915.8
$ws.Range("H12").Value = 1573.1305
$ws.Range("I12").Value = 1009.35
$ws.Range("K12").Value = 1009.35
$ws.Range("M12").Value = -839.35
$ws.Range("H40").Value = 7277.778
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825
$ws.Range("H70").Value = 2064.4614
$ws.Range("J70").Value = 2081.6667
$ws.Range("L70").Value = 6245.000100000001
$ws.Range("N70").Value = -6785.000100000001
$ws.Range("H73").Value = 2064.4614
$ws.Range("J73").Value = 2081.6667
$ws.Range("L73").Value = 6245.000100000001
$ws.Range("N73").Value = -8117.000100000001
$ws.Range("H88").Value = 2182.375
$ws.Range("J88").Value = 2549
$ws.Range("L88").Value = 2549
$ws.Range("N88").Value = -3361
$ws.Range("H91").Value = 2182.375
$ws.Range("J91").Value = 2549
$ws.Range("L91").Value = 2549
$ws.Range("N91").Value = -5357
$ws.Range("H94").Value = 14996.667
$ws.Range("I94").Value = 4995
$ws.Range("J94").Value = 35000
$ws.Range("K94").Value = 4995
$ws.Range("L94").Value = 35000
$ws.Range("M94").Value = -4544
$ws.Range("N94").Value = -35902
$ws.Range("H96").Value = 1825.8572
$ws.Range("I96").Value = 356.4
$ws.Range("J96").Value = 5499.5
$ws.Range("K96").Value = 1069.2
$ws.Range("L96").Value = 16498.5
$ws.Range("M96").Value = 303.8000000000002
$ws.Range("N96").Value = -19244.5
$ws.Range("H103").Value = 4780.5625
$ws.Range("I103").Value = 2916.6667
$ws.Range("J103").Value = 5898.9
$ws.Range("K103").Value = 8750.000100000001
$ws.Range("L103").Value = 17696.7
$ws.Range("M103").Value = -8164.000100000001
$ws.Range("N103").Value = -18868.7
$ws.Range("H107").Value = 175.15384
$ws.Range("I107").Value = 188.5
$ws.Range("J107").Value = 15
$ws.Range("K107").Value = 188.5
$ws.Range("L107").Value = 15
$ws.Range("M107").Value = 1731.5
$ws.Range("N107").Value = -3855
$ws.Range("H131").Value = 346.33334
$ws.Range("I131").Value = 346.33334
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1039.00002
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = 4000.99998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3575.0588
$ws.Range("I32").Value = 3575.0588
$ws.Range("K32").Value = 3575.0588
$ws.Range("M32").Value = -3288.0588
$ws.Range("H45").Value = 6333.1665
$ws.Range("I45").Value = 5666.3335
$ws.Range("K45").Value = 5666.3335
$ws.Range("M45").Value = -5289.3335
$ws.Range("H74").Value = 1401.3334
$ws.Range("I74").Value = 826.625
$ws.Range("K74").Value = 826.625
$ws.Range("M74").Value = 47.375
$ws.Range("H77").Value = 1401.3334
$ws.Range("I77").Value = 826.625
$ws.Range("K77").Value = 4133.125
$ws.Range("M77").Value = 234.875
$ws.Range("H132").Value = 1659.5
$ws.Range("I132").Value = 1442.2
$ws.Range("K132").Value = 4326.6
$ws.Range("M132").Value = -1796.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3104.7693
$ws.Range("J20").Value = 5139.8
$ws.Range("L20").Value = 5139.8
$ws.Range("N20").Value = -5633.8
$ws.Range("H94").Value = 1410.6
$ws.Range("J94").Value = 1517.125
$ws.Range("L94").Value = 1517.125
$ws.Range("N94").Value = -2419.125
$ws.Range("H105").Value = 4318.5
$ws.Range("I105").Value = 4242.778
$ws.Range("K105").Value = 4242.778
$ws.Range("M105").Value = -2495.778
$ws.Range("H107").Value = 683.3333
$ws.Range("I107").Value = 683.3333
$ws.Range("K107").Value = 683.3333
$ws.Range("M107").Value = 1236.6667
$ws.Range("H134").Value = 2529.1538
$ws.Range("I134").Value = 1487.9
$ws.Range("K134").Value = 4463.700000000001
$ws.Range("M134").Value = -1928.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1862.9166
$ws.Range("I7").Value = 1470
$ws.Range("J7").Value = 2648.75
$ws.Range("K7").Value = 1470
$ws.Range("L7").Value = 2648.75
$ws.Range("M7").Value = -1357
$ws.Range("N7").Value = -2874.75
$ws.Range("H58").Value = 3833.6667
$ws.Range("I58").Value = 3833.6667
$ws.Range("K58").Value = 3833.6667
$ws.Range("M58").Value = -3630.6667
$ws.Range("H62").Value = 2131.6667
$ws.Range("J62").Value = 2148
$ws.Range("L62").Value = 2148
$ws.Range("N62").Value = -3396
$ws.Range("H65").Value = 2131.6667
$ws.Range("J65").Value = 2148
$ws.Range("L65").Value = 10740
$ws.Range("N65").Value = -16980
$ws.Range("H107").Value = 452.04544
$ws.Range("I107").Value = 330.5
$ws.Range("K107").Value = 330.5
$ws.Range("M107").Value = 1589.5
$ws.Range("H136").Value = 3833.6667
$ws.Range("I136").Value = 3833.6667
$ws.Range("K136").Value = 11501.0001
$ws.Range("M136").Value = -8951.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 787.5
$ws.Range("J4").Value = 950
$ws.Range("L4").Value = 2850
$ws.Range("N4").Value = -3074
$ws.Range("H68").Value = 1049.8
$ws.Range("J68").Value = 999.6667
$ws.Range("L68").Value = 2999.0001
$ws.Range("N68").Value = -4621.0001
$ws.Range("H71").Value = 1049.8
$ws.Range("J71").Value = 999.6667
$ws.Range("L71").Value = 8997.0003
$ws.Range("N71").Value = -17109.0003
$ws.Range("H108").Value = 1100
$ws.Range("I108").Value = 1100
$ws.Range("K108").Value = 3300
$ws.Range("M108").Value = -420
$ws.Range("H120").Value = 3000
$ws.Range("I120").Value = 3000
$ws.Range("K120").Value = 9000
$ws.Range("M120").Value = -4162
$ws.Range("H131").Value = 1022
$ws.Range("J131").Value = 1042.5714
$ws.Range("L131").Value = 3127.7142
$ws.Range("N131").Value = -13207.7142
$ws.Range("H132").Value = 1426.25
$ws.Range("J132").Value = 1426.25
$ws.Range("L132").Value = 12836.25
$ws.Range("N132").Value = -17896.25
$ws.Range("H141").Value = 11953.667
$ws.Range("I141").Value = 7999
$ws.Range("K141").Value = 23997
$ws.Range("M141").Value = -18817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9843.727999999999
$ws.Range("J46").Value = 15146.833
$ws.Range("L46").Value = 15146.833
$ws.Range("N46").Value = -15458.833
$ws.Range("H123").Value = 39999
$ws.Range("J123").Value = 39999
$ws.Range("L123").Value = 39999
$ws.Range("N123").Value = -44899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2233.4614
$ws.Range("I22").Value = 2170
$ws.Range("J22").Value = 2995
$ws.Range("K22").Value = 2170
$ws.Range("L22").Value = 2995
$ws.Range("M22").Value = -1875
$ws.Range("N22").Value = -3585
$ws.Range("H27").Value = 2233.4614
$ws.Range("I27").Value = 2170
$ws.Range("J27").Value = 2995
$ws.Range("K27").Value = 2170
$ws.Range("L27").Value = 2995
$ws.Range("M27").Value = -2063
$ws.Range("N27").Value = -3209
$ws.Range("H46").Value = 2538.3
$ws.Range("I46").Value = 2360.25
$ws.Range("J46").Value = 2657
$ws.Range("K46").Value = 2360.25
$ws.Range("L46").Value = 2657
$ws.Range("M46").Value = -2172.25
$ws.Range("N46").Value = -3033

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5956.8335
$ws.Range("I62").Value = 5956.8335
$ws.Range("K62").Value = 5956.8335
$ws.Range("M62").Value = -5332.8335
$ws.Range("H65").Value = 5956.8335
$ws.Range("I65").Value = 5956.8335
$ws.Range("K65").Value = 29784.1675
$ws.Range("M65").Value = -26664.1675
$ws.Range("H126").Value = 1150.3334
$ws.Range("I126").Value = 1109.8
$ws.Range("J126").Value = 1201
$ws.Range("K126").Value = 3329.4
$ws.Range("L126").Value = 3603
$ws.Range("M126").Value = -859.3999999999996
$ws.Range("N126").Value = -8543
$ws.Range("H129").Value = 80000
$ws.Range("J129").Value = 80000
$ws.Range("L129").Value = 80000
$ws.Range("N129").Value = -90000
$ws.Range("H132").Value = 167757.5
$ws.Range("I132").Value = 201109
$ws.Range("K132").Value = 603327
$ws.Range("M132").Value = -600797
